# Update the cryptocurrency price table with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) contains values such as "69.625.12", "0.999", "1.00"
# that must stay as literal text (Excel would otherwise coerce them into
# numbers and silently drop meaningful trailing/duplicate digits). Force the
# whole column to text formatting before writing any of the new values.
$ws.Range("D2:D51").NumberFormat = "@"

# NOTE: this PowerShell-like engine only supports positional parameter
# binding for user functions, so Set-Row takes its arguments positionally:
#   Set-Row <row> <coin> <link> <price> <volume>
# Pass $null for any column that should be left untouched.
function Set-Row($Row, $Coin, $Link, $Price, $Volume) {
    if ($Coin -ne $null) {
        $ws.Cells.Item($Row, 2).Value = $Coin
    }
    if ($Link -ne $null) {
        $ws.Cells.Item($Row, 3).Value = $Link
    }
    if ($Price -ne $null) {
        $ws.Cells.Item($Row, 4).Value = $Price
    }
    if ($Volume -ne $null) {
        $ws.Cells.Item($Row, 5).Value = $Volume
    }
}

Set-Row 2  $null $null "69.612.69"  "  +0.30%  "
Set-Row 3  $null $null "3.490.51"   "  +0.10%  "
Set-Row 4  $null $null "0.997"      "  -0.12%  "
Set-Row 5  $null $null "603.20"     "  -1.05%  "
Set-Row 6  $null $null "193.86"     "  +4.19%  "
Set-Row 7  $null $null $null        "  -0.09%  "
Set-Row 8  $null $null $null        "  -0.01%  "
Set-Row 9  $null $null "0.200"      "  -6.69%  "
Set-Row 10 $null $null "0.647"      "  +0.09%  "
Set-Row 11 $null $null "53.15"      "  +0.38%  "
Set-Row 12 $null $null $null        "  -2.61%  "
Set-Row 14 $null $null "4.067.66"   "  +0.70%  "
Set-Row 15 $null $null "593.22"     "  -1.52%  "
Set-Row 16 $null $null "69.746.55"  "  +0.42%  "
Set-Row 17 $null $null "12.74"      "  +1.62%  "
Set-Row 18 $null $null "18.93"      "  +0.62%  "
Set-Row 19 $null $null "3.503.90"   "  +0.27%  "
Set-Row 21 $null $null $null        "  -0.15%  "
Set-Row 22 $null $null "18.05"      "  +5.39%  "
Set-Row 23 $null $null "5.31"       "  +4.25%  "
Set-Row 24 $null $null $null        "  +0.62%  "
Set-Row 25 $null $null "101.92"     "  -3.82%  "
Set-Row 26 $null $null "3.14"       "  +3.04%  "
Set-Row 27 $null $null "10.83"      "  -1.05%  "
Set-Row 28 $null $null "9.52"       "  -1.62%  "
Set-Row 29 $null $null "33.22"      "  -0.87%  "
Set-Row 30 $null $null $null        "  +0.63%  "
Set-Row 31 $null $null "4.17"       "  +1.98%  "
Set-Row 32 $null $null $null        "  -0.33%  "
Set-Row 33 $null $null "0.114"      "  -0.35%  "
Set-Row 34 $null $null "63.09"      "  -0.29%  "
Set-Row 35 $null $null "0.0₃0830"   "  +6.65%  "
Set-Row 36 $null $null "3.709.89"   "  +2.72%  "

# Rows 37 and 38 swap identities (Fetch.AI now ranks above Dai).
Set-Row 37 "Fetch.AI" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet" "3.08" "  -3.00%  "
Set-Row 38 "Dai"      "https://coinranking.com/coin/MoTuySvg7+dai-dai"          "1.00" "  +0.20%  "

Set-Row 39 $null $null "3.64"       "  -0.68%  "
Set-Row 40 $null $null $null        "  -1.68%  "
Set-Row 41 $null $null "36.31"      $null
Set-Row 42 $null $null "478.90"     "  -8.20%  "
Set-Row 43 $null $null $null        "  -3.21%  "
Set-Row 44 $null $null "0.0451"     "  -2.06%  "
Set-Row 45 $null $null $null        "  -1.97%  "
Set-Row 46 $null $null $null        "  -4.71%  "
Set-Row 47 $null $null "3.27"       "  -1.99%  "
Set-Row 48 $null $null "1.00"       "  +0.25%  "
Set-Row 49 $null $null "8.40"       "  -4.47%  "
Set-Row 50 $null $null "0.000243"   "  +0.61%  "
Set-Row 51 $null $null "1.28"       "  +9.86%  "
